# Update LocationApi endpoints to the new iassureit.com-hosted URLs and turn
# them into real hyperlinks (Excel auto-creates the "Hyperlink" cell style /
# underlined theme-10 font + cellStyleXf the first time Hyperlinks.Add runs
# on a workbook that doesn't have one yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the endpoint text first --------------------------------------
# (states, districts, areas/one, blocks, cities, areas/list -- this is the
# order the text was originally authored in, and it controls the order new
# entries land in the shared-string table.)
$ws.Range("C2").Value2 = "http://locationapi.iassureit.com/api/states/get/list/:countryCode"
$ws.Range("C3").Value2 = "http://locationapi.iassureit.com/api/districts/get/list/:stateCode/:countryCode"
$ws.Range("C7").Value2 = "http://locationapi.iassureit.com/api/areas/get/one/:pincode"
$ws.Range("C4").Value2 = "http://locationapi.iassureit.com/api/blocks/get/list/:districtName/:stateCode/:countryCode"
$ws.Range("C5").Value2 = "http://locationapi.iassureit.com/api/cities/get/list/:blockName/:districtName/:stateCode/:countryCode"
$ws.Range("C6").Value2 = "http://locationapi.iassureit.com/api/areas/get/list/:cityName/:blockName/:districtName/:stateCode/:countryCode"

# --- Then turn each endpoint cell into a hyperlink, in sheet (ref) order ---
$ws.Hyperlinks.Add($ws.Range("C2"), "http://locationapi.iassureit.com/api/states/get/list/:countryCode")
$ws.Hyperlinks.Add($ws.Range("C3"), "http://locationapi.iassureit.com/api/districts/get/list/:stateCode/:countryCode")
$ws.Hyperlinks.Add($ws.Range("C4"), "http://locationapi.iassureit.com/api/blocks/get/list/:districtName/:stateCode/:countryCode")
$ws.Hyperlinks.Add($ws.Range("C5"), "http://locationapi.iassureit.com/api/cities/get/list/:blockName/:districtName/:stateCode/:countryCode")
$ws.Hyperlinks.Add($ws.Range("C6"), "http://locationapi.iassureit.com/api/areas/get/list/:cityName/:blockName/:districtName/:stateCode/:countryCode")
$ws.Hyperlinks.Add($ws.Range("C7"), "http://locationapi.iassureit.com/api/areas/get/one/:pincode")

# --- Restore the active selection (it had moved from D7 to C7) ------------
$ws.Range("C7").Select()
